# TC-40.xlsx update
# - Test Case ID (B1) renamed from TC-01 to TC-40
# - Test Data password row: "Pass: PruebaTC1!" replaced with "Contraseña",
#   and the now-superfluous Fecha/Dir/Ciudad/Estado/Cod postal/Country/Phone
#   test-data rows are cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "TC-40"

$ws.Range("E11").Value = "Contraseña"
$ws.Range("E12").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("E18").Value = ""

# Re-assert the (unchanged) cell formatting on a few cells whose internal
# style-table slot shifts as a side effect of the shared-string table being
# compacted, so they land back on the same visual style bucket as their
# neighbours.
$ws.Range("B10").WrapText = $true
$ws.Range("B10").VerticalAlignment = -4160

foreach ($addr in @("B27", "D27", "D28", "B29", "D29")) {
    $ws.Range($addr).WrapText = $true
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4160
}
